# Edit: insert two new price records at rows 164-165 (pushing the
# existing rows 164-239 down to 166-241), matching a weekly data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 164; this pushes the old rows
# 164-239 down to 166-241 and grows the sheet's dimension automatically.
$ws.Rows("164:165").Insert()

# --- New row 164 ---
$ws.Range("A164").Value2 = 1
$ws.Range("B164").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C164").Value2 = "Arica y Parinacota"
$ws.Range("D164").Value2 = 45119
$ws.Range("E164").Value2 = 15
$ws.Range("F164").Value2 = 100114001
$ws.Range("G164").Value2 = "Papa"
$ws.Range("H164").Value2 = "Cardinal"
$ws.Range("I164").Value2 = "1a (cosecha)"
$ws.Range("J164").Value2 = 1000
$ws.Range("K164").Value2 = 16000
$ws.Range("L164").Value2 = 17000
$ws.Range("M164").Value2 = 16500
$ws.Range("N164").Value2 = "`$/saco 25 kilos"
$ws.Range("O164").Value2 = "Provincia de Melipilla"
$ws.Range("P164").Value2 = 660
$ws.Range("Q164").Value2 = 25
$ws.Range("R164").Value2 = "Hortaliza"

# --- New row 165 ---
$ws.Range("A165").Value2 = 1
$ws.Range("B165").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C165").Value2 = "Arica y Parinacota"
$ws.Range("D165").Value2 = 45119
$ws.Range("E165").Value2 = 15
$ws.Range("F165").Value2 = 100114001
$ws.Range("G165").Value2 = "Papa"
$ws.Range("H165").Value2 = "Rodeo"
$ws.Range("I165").Value2 = "1a (guarda)"
$ws.Range("J165").Value2 = 1000
$ws.Range("K165").Value2 = 14000
$ws.Range("L165").Value2 = 15000
$ws.Range("M165").Value2 = 14500
$ws.Range("N165").Value2 = "`$/saco 25 kilos"
$ws.Range("O165").Value2 = "Región de La Araucanía"
$ws.Range("P165").Value2 = 580
$ws.Range("Q165").Value2 = 25
$ws.Range("R165").Value2 = "Hortaliza"
